$d = $word.ActiveDocument

$lines = @(
    'cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565',
    '6/23/23, 5:37 PM',
    'CITY PLANNING COMMISSION',
    'MEETING AGENDA',
    'TUESDAY, JUNE 13, 2023',
    '01:30 P.M. COUNCIL CHAMBERS',
    'Minutes',
    '1. Adoption of the minutes of the May 23, 2023 meeting Minutes for the May 23',
    'meeting',
    'Business',
    '2. Zoning Docket 039/23 - Deferred from May 23, 2023 Applicant(s): Bywater JEJ,',
    'LLC Request: Conditional use to permit a hotel in an HU-MU Historic',
    'Marigny/Tremé/Bywater Mixed-Use District Property description: Square 292, Lots',
    'Y, 10, 23, and 24 or Lots Y, 1, and two undesignated lots, in the Third Municipal',
    'District, bounded by Bartholomew Street, Burgundy Street, Mazant Street, and North',
    'Rampart Street Address(es): 901 Bartholemew Street and 4015 Burgundy Street',
    'Preliminary staff report',
    'Zoning Docket 042/23 Applicant(s): Professional Funeral Services, Inc. Request:',
    '3.',
    'Zoning change from an LI Light Industrial District to an MU-1 Medium Intensity',
    'Mixed-Use District Property description: The entirety of Square 1804, in the Third',
    'Municipal District, bounded by Franklin Avenue, Abundance Street, Agriculture',
    'Street, and Eads Street Municipal address(es): 2905-2941 Franklin Avenue',
    'Preliminary staff report',
    'Zoning Docket 043/23 Applicant(s): New Orleans Property Partners, LLC Request:',
    '4.',
    'Conditional use to permit the retail sale of packaged alcoholic beverages in a CBD-5',
    'Urban Core Neighborhood Lower Intensity Mixed-Use District and the rescission or',
    'amendment of Ordinance No. 21,252 MCS (Zoning Docket 001/03) Property',
    'description: Square 256, Parcel 2 or Lots A and 37, in the First Municipal District,',
    'bounded by Howard Avenue, Baronne Street, O''Keefe Avenue, and Julia Street',
    'Municipal address(es): 929 Howard Avenue',
    'Preliminary staff report',
    'Zoning Docket 044/23 Applicant(s): City Council Motion No. M-23-146 Request:',
    '5.',
    'Amendment to the text of the Comprehensive Zoning Ordinance to modify Article 24,',
    'Section 24.15.C to read as follows: "24.15.C MAINTENANCE AND RELOCATION',
    'A designated classic sign may be moved to a new location when and where the',
    'associated business moves, or an exact replica of a designated classic sign may be',
    'installed at a new location when a business with a classic sign moves. Only the',
    'business name, other descriptive language, and color of an original designated classic',
    'sign may be changed to reflect a new business at the original classic sign location."',
    'https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565',
    '1/3',
    'cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565',
    '6/23/23, 5:37 PM',
    'Preliminary Staff Report',
    'Zoning Docket 045/23 Applicant(s): City Council Motion M-23-162 Request:',
    '6.',
    'Amendment to the text of the Comprehensive Zoning Ordinance to modify Article 18,',
    'Section 18.30, the University Area Off-Street Parking Overlay District, the intent of',
    'which is to require off-street vehicular parking for certain increases in the number of',
    'existing bedrooms and bathrooms, (i) to expand the area of applicability to include',
    'additional neighborhoods, including parts of neighborhoods bounded by Front Street,',
    'Jefferson Avenue, South Claiborne Avenue, Toledano Street, Earhart Boulevard and',
    'South Carrollton Avenue; (ii) redefine "bedroom" to include a window or a door that',
    'opens to a space other than another sleeping room, bathroom or toilet room, and',
    'require that every shared bedroom be at least 50 sq. ft. for each occupant; (iii) lower',
    'the threshold for off-street parking requirement for new construction or renovations',
    'from four bedrooms and three bathrooms to three bedrooms and two bathrooms; (iv)',
    'redefine "full bathroom" to include a bathtub and/or shower; (v) require new',
    'construction or renovation to include a common room which cannot be a bedroom;',
    '(vi) require that new parking spaces be a minimum of 8 1/2 feet wide by 18 feet long;',
    '(vii) prohibit tandem parking; (viii) recognize possible shared parking arrangements',
    'between adjacent lots; and (viv) prohibit the raising of structures to accommodate off-',
    'street parking.',
    'Preliminary staff report',
    '7. Zoning Docket 046/23 Applicant(s): City Council Motion M-23-M Request:',
    'Amendment to the text of the Comprehensive Zoning Ordinance to modify Article 1,',
    'Section 17.3.A to classify "broadcasting studio" as a permitted use in the CBD-5',
    'Urban Core Neighborhood Lower Intensity Mixed-Use District.',
    'Preliminary staff report',
    '8. Firearm Sales Study The purpose of this study shall be to make recommendations',
    'regarding the addition of "firearm sales" as a use separate and distinct from "retail',
    'goods establishment" and to make recommendations regarding potential zoning',
    'changes that would prohibit the sale of firearms within 1,000 feet of educational',
    'facilities.',
    'Preliminary Staff Report',
    'Subdivision Ratifications This ratification includes applications that have not yet been',
    '9.',
    'certified by the staff and where the Commission is ratifying the staff determination',
    'that the application may be reviewed and approved administratively.',
    'List of subdivisions deferred from the May 23 meeting',
    'List of new subdivisions for the June 13 meeting',
    '2. Zoning Docket 039/23 - Deferred from May 23, 2023 Applicant(s): Bywater JEJ, LLC',
    'Request: Conditional use to permit a hotel in an HU-MU Historic',
    'https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565',
    '2/3',
    'cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565',
    '6/23/23, 5:37 PM',
    'Marigny/Tremé/Bywater Mixed-Use District Property description: Square 292, Lots Y,',
    '10, 23, and 24 or Lots Y, 1, and two undesignated lots, in the Third Municipal District,',
    'bounded by Bartholomew Street, Burgundy Street, Mazant Street, and North Rampart',
    'Street Address(es): 901 Bartholemew Street and 4015 Burgundy Street',
    'Zoning Docket 042/23 Applicant(s): Professional Funeral Services, Inc. Request:',
    '3.',
    'Zoning change from an LI Light Industrial District to an MU-1 Medium Intensity',
    'Mixed-Use District Property description: The entirety of Square 1804, in the Third',
    'Municipal District, bounded by Franklin Avenue, Abundance Street, Agriculture',
    'Street, and Eads Street Municipal address(es): 2905-2941 Franklin Avenue',
    'Zoning Docket 042/23 Applicant(s): Professional Funeral Services, Inc. Request:',
    '3.',
    'Zoning change from an LI Light Industrial District to an MU-1 Medium Intensity',
    'Mixed-Use District Property description: The entirety of Square 1804, in the Third',
    'Municipal District, bounded by Franklin Avenue, Abundance Street, Agriculture',
    'Street, and Eads Street Municipal address(es): 2905-2941 Franklin Avenue',
    'Robert D. Rivers',
    'Executive Director',
    '3/3',
    'https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4565'
)

$vt = [char]11
$cr = [char]13

$bodyText = ($lines -join $vt)
$para2Text = "" + $cr + $vt + '---------- End of Page 1 ----------' + $vt

$fullText = $bodyText + $para2Text

$r = $d.Paragraphs.First.Range
$r.Text = $fullText

if ($d.Paragraphs.Count -ne 2) {
    throw ("Expected 2 paragraphs after edit, got " + $d.Paragraphs.Count)
}

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
Write-Output ("P1 starts: " + $d.Paragraphs.Item(1).Range.Text.Substring(0, 40))
Write-Output ("P2: " + $d.Paragraphs.Item(2).Range.Text)
